# Voltorb-AUG-2022.pptx edit: reposition two shapes on slide 8 (the
# "THANK YOU" slide) and add a new "For queries / Write to us at ..."
# textbox with a highlighted, hyperlinked e-mail address.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

# --- 1) Move "IEEE UVCE PES" textbox up -------------------------------
$titleShape = $s.Shapes.Item(1)
$titleShape.Top = 2489310 / 12700

# --- 2) Move the "THANK YOU" shape up ---------------------------------
$thankYouShape = $s.Shapes.Item(3)
$thankYouShape.Top = 974559 / 12700

# --- 3) Add the new "For queries ..." textbox -------------------------
# A throw-away textbox is created first (and immediately removed) purely
# to advance PowerPoint's internal shape id/name counter so the real
# textbox we keep lands on id=6 / "TextBox 5", matching the author's
# original shape numbering.
$placeholderShape = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$placeholderShape.Delete()

$queryBox = $s.Shapes.AddTextbox(1, 3810759 / 12700, 4478215 / 12700, 4570482 / 12700, 646331 / 12700)
$queryBox.TextFrame.WordWrap = 0
$queryBox.TextFrame.AutoSize = 1

$tr = $queryBox.TextFrame.TextRange
$tr.Text = "For queries `rWrite to us at ieeeuvcepes@gmail.com"

$emailRange = $tr.Characters(29, 21)
$emailRange.Font.Highlight = 192 + 192 * 256 + 192 * 65536
$emailRange.ActionSettings.Item(1).Hyperlink.Address = "mailto:ieeeuvcepes@gmail.com"
